$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New word pairs (misspelled -> correct) for rows 16-21, columns K (11) and L (12)
$pairs = @(
    @("pleague", "plague"),
    @("autum", "autumn"),
    @("Syberia", "Siberia"),
    @("anticeptive", "antiseptic"),
    @("nurf", "nerf"),
    @("obeisity", "obesity")
)

$row = 16
foreach ($pair in $pairs) {
    $ws.Cells.Item($row, 11).Value = $pair[0]
    $ws.Cells.Item($row, 12).Value = $pair[1]
    $row++
}

# Update the active selection to match the final state
$ws.Range("L22").Select()
